# Auto commit at 2025-12-26  9:19:37.04
#
# Updates the "Metrics" sheet's monthly/annual KPI values, which ripple
# through the "today" sheet's formulas, and switches the active sheet
# selection from "Chargingdata" to "today".

$wb = $excel.ActiveWorkbook

# --- Update the metric values on the "Metrics" sheet ------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 342173.67000000004
$metrics.Range("B3").Value = 293204.55999999994
$metrics.Range("B4").Value = 104649.23999999999
$metrics.Range("B5").Value = 13940
$metrics.Range("B6").Value = 5544880.7800000003
$metrics.Range("B7").Value = 4693557.5200000005
$metrics.Range("B8").Value = 1636606.12
$metrics.Range("B9").Value = 216647
$metrics.Range("B10").Value = 34010261.769999996
$metrics.Range("B11").Value = 31968832.68
$metrics.Range("B12").Value = 11918328.159999995
$metrics.Range("B13").Value = 1314277

# The selection on the Metrics sheet moved from C21 to D13.
$metrics.Range("D13").Select()

# --- Move the active tab from "Chargingdata" to "today" ---------------
$chargingData = $wb.Worksheets.Item("Chargingdata")
$today = $wb.Worksheets.Item("today")

$today.Activate()

# The selection on the "today" sheet moved from F6 to H14.
$today.Range("H14").Select()
